$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.719.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.44%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.695.91'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.02%  '
# Row 4
$ws.Range("E4").Value = '  +0.34%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.99%  '
# Row 6
$ws.Range("E6").Value = '  +0.40%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3957'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4074'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.80%  '
# Row 9
$ws.Range("E9").Value = '  -2.09%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.005'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.41%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.96'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.81%  '
# Row 12
$ws.Range("E12").Value = '  +1.92%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.282'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.37%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.70%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.041'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.56%  '
# Row 16
$ws.Range("E16").Value = '  +0.51%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.698.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.20%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '100.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.43%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07039'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.25%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.93%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.995'
$ws.Range("D21").Style = "Normal"
# Row 22
$ws.Range("E22").Value = '  +0.22%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.17%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.693.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.36%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.279'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.09%  '
# Row 26
$ws.Range("E26").Value = '  +2.43%  '
# Row 27
$ws.Range("E27").Value = '  +1.35%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.65%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '136.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.18%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.171'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.28%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.506'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.08%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08699'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.54%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.054'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.98%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.077'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.72%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.93%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2743'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.89%  '
# Row 37
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '14.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.75%  '
# Row 38
$ws.Range("B38").Value = 'WEMIXTOKEN'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.881'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.11%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09251'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.07%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02728'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.32%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.476'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.18%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7671'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.40%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.26'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.63%  '
# Row 44
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7181'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.01%  '
# Row 45
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.591'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.27%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.220'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.29%  '
# Row 47
$ws.Range("E47").Value = '  +0.40%  '
# Row 48
$ws.Range("E48").Value = '  +0.00%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.323'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.46%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '91.11'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.25%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07978'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.34%  '
